# Weekly data refresh: a new week's price record is inserted at the top of
# the data block (row 49), pushing all subsequent records down by one row
# (old row 49 -> new row 50, ..., old row 69 -> new row 70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49; this shifts rows 49:69 down to 50:70 and carries
# the existing row formatting (e.g. the date style on column D) along.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with this week's record.
$ws.Cells.Item(49, 1).Value = 4
$ws.Cells.Item(49, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(49, 3).Value = "Los Lagos"
$ws.Cells.Item(49, 4).Value = 45097
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = 100112012
$ws.Cells.Item(49, 7).Value = "Espinaca"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 25
$ws.Cells.Item(49, 11).Value = 13000
$ws.Cells.Item(49, 12).Value = 13000
$ws.Cells.Item(49, 13).Value = 13000
$ws.Cells.Item(49, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 1300
$ws.Cells.Item(49, 17).Value = 10
$ws.Cells.Item(49, 18).Value = "Hortaliza"
